$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds data cells per the 2025-11-19 refresh

# Row 2
$ws.Range("F2").Value = 3.4
$ws.Range("G2").Value = 3.9
$ws.Range("H2").Value = 2.5
$ws.Range("I2").Value = 2.76
$ws.Range("J2").Value = 2.76
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1.58
$ws.Range("M2").Value = 1.12
$ws.Range("N2").Value = 2.72
$ws.Range("P2").Value = 1.53
$ws.Range("Q2").Value = 2.62
$ws.Range("S2").Value = 5.1
$ws.Range("T2").Value = 2.02
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 1.35
$ws.Range("X2").Value = 8.800000000000001
$ws.Range("Z2").Value = 16
$ws.Range("AA2").Value = 44
$ws.Range("AB2").Value = 10.5
$ws.Range("AC2").Value = 7
$ws.Range("AD2").Value = 13
$ws.Range("AF2").Value = 25
$ws.Range("AH2").Value = 21
$ws.Range("AL2").Value = 300
$ws.Range("AO2").Value = 600

# Row 3
$ws.Range("F3").Value = 3.25
$ws.Range("G3").Value = 3.65
$ws.Range("H3").Value = 2.38
$ws.Range("I3").Value = 2.6
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 3.3
$ws.Range("L3").Value = 1.52
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 1.43
$ws.Range("P3").Value = 1.65
$ws.Range("Q3").Value = 2.32
$ws.Range("R3").Value = 1.24
$ws.Range("S3").Value = 4.5
$ws.Range("T3").Value = 1.93
$ws.Range("U3").Value = 1.9
$ws.Range("V3").Value = 1.64
$ws.Range("W3").Value = 1.38
$ws.Range("X3").Value = 11
$ws.Range("Y3").Value = 8.800000000000001
$ws.Range("Z3").Value = 16
$ws.Range("AB3").Value = 11
$ws.Range("AD3").Value = 12
$ws.Range("AE3").Value = 110
$ws.Range("AF3").Value = 24
$ws.Range("AG3").Value = 15
$ws.Range("AH3").Value = 21
$ws.Range("AI3").Value = 100
$ws.Range("AJ3").Value = 170
$ws.Range("AK3").Value = 75
$ws.Range("AL3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 44

# Row 4
$ws.Range("F4").Value = 1.32
$ws.Range("G4").Value = 1.33
$ws.Range("H4").Value = 12
$ws.Range("I4").Value = 12.5
$ws.Range("K4").Value = 6.2
$ws.Range("L4").Value = 1.36
$ws.Range("N4").Value = 4.4
$ws.Range("O4").Value = 1.28
$ws.Range("P4").Value = 2.14
$ws.Range("Q4").Value = 1.84
$ws.Range("R4").Value = 1.44
$ws.Range("S4").Value = 3.15
$ws.Range("T4").Value = 2.4
$ws.Range("U4").Value = 1.69
$ws.Range("X4").Value = 18.5
$ws.Range("Y4").Value = 36
$ws.Range("Z4").Value = 120
$ws.Range("AA4").Value = 620
$ws.Range("AB4").Value = 7.6
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 48
$ws.Range("AE4").Value = 310
$ws.Range("AF4").Value = 7
$ws.Range("AH4").Value = 38
$ws.Range("AI4").Value = 210
$ws.Range("AK4").Value = 15.5
$ws.Range("AL4").Value = 46
$ws.Range("AM4").Value = 240
$ws.Range("AN4").Value = 5.7
$ws.Range("AO4").Value = 360

# Row 5
$ws.Range("G5").Value = 1.91
$ws.Range("H5").Value = 5.1
$ws.Range("I5").Value = 5.4
$ws.Range("J5").Value = 3.55
$ws.Range("K5").Value = 3.7
$ws.Range("L5").Value = 1.57
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 2.78
$ws.Range("O5").Value = 1.52
$ws.Range("R5").Value = 1.2
$ws.Range("S5").Value = 5.3
$ws.Range("T5").Value = 2.22
$ws.Range("U5").Value = 1.69
$ws.Range("V5").Value = 1.22
$ws.Range("W5").Value = 2.08
$ws.Range("X5").Value = 10
$ws.Range("Y5").Value = 14
$ws.Range("Z5").Value = 38
$ws.Range("AA5").Value = 170
$ws.Range("AB5").Value = 6.8
$ws.Range("AC5").Value = 8.199999999999999
$ws.Range("AD5").Value = 23
$ws.Range("AE5").Value = 100
$ws.Range("AF5").Value = 9.800000000000001
$ws.Range("AG5").Value = 11
$ws.Range("AH5").Value = 28
$ws.Range("AI5").Value = 130
$ws.Range("AJ5").Value = 22
$ws.Range("AK5").Value = 25
$ws.Range("AL5").Value = 60
$ws.Range("AM5").Value = 220
$ws.Range("AN5").Value = 21
$ws.Range("AO5").Value = 160

# Row 6
$ws.Range("G6").Value = 2.82
$ws.Range("H6").Value = 2.8
$ws.Range("I6").Value = 3.1
$ws.Range("J6").Value = 3.25
$ws.Range("K6").Value = 3.65
$ws.Range("L6").Value = 1.45
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 3.4
$ws.Range("P6").Value = 1.85
$ws.Range("Q6").Value = 2.02
$ws.Range("R6").Value = 1.32
$ws.Range("S6").Value = 3.5
$ws.Range("T6").Value = 1.72
$ws.Range("U6").Value = 2.08
$ws.Range("V6").Value = 1.48
$ws.Range("W6").Value = 1.54
$ws.Range("AA6").Value = 150
$ws.Range("AE6").Value = 90
$ws.Range("AH6").Value = 42
$ws.Range("AI6").Value = 260
$ws.Range("AJ6").Value = 170
$ws.Range("AK6").Value = 80

# Row 7
$ws.Range("F7").Value = 1.93
$ws.Range("G7").Value = 2.02
$ws.Range("I7").Value = 5.2
$ws.Range("J7").Value = 3.3
$ws.Range("K7").Value = 3.55
$ws.Range("L7").Value = 1.52
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 2.96
$ws.Range("O7").Value = 1.46
$ws.Range("P7").Value = 1.65
$ws.Range("Q7").Value = 2.38
$ws.Range("R7").Value = 1.23
$ws.Range("S7").Value = 4.8
$ws.Range("T7").Value = 2.06
$ws.Range("U7").Value = 1.79
$ws.Range("V7").Value = 1.24
$ws.Range("W7").Value = 1.99
$ws.Range("X7").Value = 14
$ws.Range("Y7").Value = 28
$ws.Range("AA7").Value = 150
$ws.Range("AB7").Value = 7.2
$ws.Range("AC7").Value = 11.5
$ws.Range("AD7").Value = 46
$ws.Range("AF7").Value = 21
$ws.Range("AG7").Value = 11
$ws.Range("AH7").Value = 85
$ws.Range("AJ7").Value = 28
$ws.Range("AK7").Value = 70
$ws.Range("AL7").Value = 150
$ws.Range("AM7").Value = 200
$ws.Range("AN7").Value = 22

# Row 8
$ws.Range("F8").Value = 3.95
$ws.Range("G8").Value = 4.1
$ws.Range("H8").Value = 2.24
$ws.Range("I8").Value = 2.26
$ws.Range("J8").Value = 3.2
$ws.Range("K8").Value = 3.3
$ws.Range("L8").Value = 1.61
$ws.Range("M8").Value = 1.14
$ws.Range("N8").Value = 2.56
$ws.Range("O8").Value = 1.62
$ws.Range("P8").Value = 1.52
$ws.Range("Q8").Value = 2.82
$ws.Range("R8").Value = 1.18
$ws.Range("S8").Value = 6.2
$ws.Range("T8").Value = 2.46
$ws.Range("U8").Value = 1.65
$ws.Range("V8").Value = 1.79
$ws.Range("W8").Value = 1.32
$ws.Range("X8").Value = 8
$ws.Range("Z8").Value = 11.5
$ws.Range("AA8").Value = 27
$ws.Range("AB8").Value = 9.800000000000001
$ws.Range("AC8").Value = 7.8
$ws.Range("AD8").Value = 13
$ws.Range("AE8").Value = 34
$ws.Range("AF8").Value = 25
$ws.Range("AH8").Value = 27
$ws.Range("AI8").Value = 70
$ws.Range("AJ8").Value = 95
$ws.Range("AL8").Value = 110
$ws.Range("AM8").Value = 230
$ws.Range("AN8").Value = 980
$ws.Range("AO8").Value = 34

# Row 9
$ws.Range("F9").Value = 2.04
$ws.Range("G9").Value = 2.06
$ws.Range("H9").Value = 4.1
$ws.Range("I9").Value = 4.2
$ws.Range("J9").Value = 3.7
$ws.Range("K9").Value = 3.75
$ws.Range("L9").Value = 1.41
$ws.Range("N9").Value = 3.95
$ws.Range("O9").Value = 1.31
$ws.Range("P9").Value = 2.02
$ws.Range("Q9").Value = 1.94
$ws.Range("R9").Value = 1.39
$ws.Range("S9").Value = 3.45
$ws.Range("T9").Value = 1.81
$ws.Range("U9").Value = 2.18
$ws.Range("V9").Value = 1.31
$ws.Range("W9").Value = 1.94
$ws.Range("X9").Value = 15
$ws.Range("Y9").Value = 15.5
$ws.Range("Z9").Value = 29
$ws.Range("AA9").Value = 85
$ws.Range("AC9").Value = 8.199999999999999
$ws.Range("AD9").Value = 16
$ws.Range("AE9").Value = 48
$ws.Range("AF9").Value = 12
$ws.Range("AH9").Value = 18
$ws.Range("AI9").Value = 55
$ws.Range("AJ9").Value = 23
$ws.Range("AK9").Value = 20
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 95
$ws.Range("AN9").Value = 15.5
$ws.Range("AO9").Value = 50

# Row 10
$ws.Range("F10").Value = 2.2
$ws.Range("G10").Value = 2.24
$ws.Range("H10").Value = 3.6
$ws.Range("I10").Value = 3.7
$ws.Range("J10").Value = 3.6
$ws.Range("K10").Value = 3.7
$ws.Range("L10").Value = 1.42
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 3.85
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 1.98
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = 1.37
$ws.Range("S10").Value = 3.6
$ws.Range("T10").Value = 1.81
$ws.Range("U10").Value = 2.16
$ws.Range("V10").Value = 1.37
$ws.Range("W10").Value = 1.81
$ws.Range("X10").Value = 14.5
$ws.Range("Y10").Value = 14
$ws.Range("Z10").Value = 26
$ws.Range("AA10").Value = 70
$ws.Range("AB10").Value = 9.800000000000001
$ws.Range("AC10").Value = 7.8
$ws.Range("AD10").Value = 15
$ws.Range("AF10").Value = 13
$ws.Range("AG10").Value = 10.5
$ws.Range("AH10").Value = 17.5
$ws.Range("AJ10").Value = 27
$ws.Range("AK10").Value = 22
$ws.Range("AL10").Value = 36
$ws.Range("AN10").Value = 16.5
$ws.Range("AO10").Value = 42

# Row 11
$ws.Range("G11").Value = 2.74
$ws.Range("H11").Value = 2.9
$ws.Range("I11").Value = 3.1
$ws.Range("J11").Value = 3.4
$ws.Range("L11").Value = 1.38
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 4.2
$ws.Range("P11").Value = 2.1
$ws.Range("Q11").Value = 1.86
$ws.Range("R11").Value = 1.44
$ws.Range("S11").Value = 3.1
$ws.Range("T11").Value = 1.63
$ws.Range("U11").Value = 2.34
$ws.Range("W11").Value = 1.58
$ws.Range("AC11").Value = 8.6
$ws.Range("AG11").Value = 13
$ws.Range("AK11").Value = 34
$ws.Range("AM11").Value = 75
$ws.Range("AN11").Value = 24
